$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "part of speech" columns (H and J) so that each bare word is
# replaced by a Python-list-style string "['Word','partOfSpeech']".
$ws.Range("H2").Value  = "['Bit','noun']"
$ws.Range("J2").Value  = "['Months','noun']"
$ws.Range("H3").Value  = "['Avoid','verb']"
$ws.Range("J3").Value  = "['Appeal','verb']"
$ws.Range("H4").Value  = "['Cover','verb']"
$ws.Range("J4").Value  = "['Promo','noun']"
$ws.Range("H5").Value  = "['Shows','verb']"
$ws.Range("J5").Value  = "['Artist','noun']"
$ws.Range("H6").Value  = "['Make','verb']"
$ws.Range("J6").Value  = "['Ecstasy','noun']"
$ws.Range("H7").Value  = "['Everyone','noun']"
$ws.Range("J7").Value  = "['ure','suffix']"
$ws.Range("H8").Value  = "['Filing','verb']"
$ws.Range("J8").Value  = "['Cause','verb']"
$ws.Range("H9").Value  = "['Black','verb']"
$ws.Range("J9").Value  = "['Beliefs','noun']"
$ws.Range("H10").Value = "['Allow','verb']"
$ws.Range("J10").Value = "['Setting','verb']"
$ws.Range("H11").Value = "['Basics','adjective']"
$ws.Range("J11").Value = "['Linen','noun']"

# Adjust column widths to match the widened/added custom column widths.
$ws.Columns.Item(2).ColumnWidth  = 15.666666666666666   # B -> 16.5
$ws.Columns.Item(3).ColumnWidth  = 8                     # C -> 8.83203125
$ws.Columns.Item(4).ColumnWidth  = 8                     # D -> 8.83203125
$ws.Columns.Item(6).ColumnWidth  = 8                     # F -> 8.83203125
$ws.Columns.Item(7).ColumnWidth  = 5                     # G -> 5.83203125
$ws.Columns.Item(8).ColumnWidth  = 13                    # H -> 13.83203125
$ws.Columns.Item(9).ColumnWidth  = 8                     # I -> 8.83203125
$ws.Columns.Item(10).ColumnWidth = 13.666666666666666    # J -> 14.5

# Update the view: zoom level and the active selection/cell.
$excel.ActiveWindow.Zoom = 171
$ws.Range("C17").Select()
